# The "Microwave: 3D Audio (Humming), auto stop" TODO item (row 2) is now
# obsolete / duplicated by the already-existing "auto stop after 30s, 3D
# audio" behaviour described elsewhere in the sheet, so the whole row is
# removed from the table. Deleting the entire row shifts every row below it
# up by one, which is exactly what the target workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Rows.Item(2).Delete()
